$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TrpLoop2b_NOE")

# New column E: <r^-6>^-1/6 (nm) header + per-row formula converting the
# r^-6-averaged distance (column B) into an <r^-6>^-1/6 distance in nm.
$ws.Range("E1").Value = "<r^-6>^-1/6 (nm)"

# Fill the formula down in chunks (matches how the workbook was actually
# authored -- separate fill/paste operations produce separate shared-formula
# groups in the saved XML).
$ws.Range("E2").Formula = "=B2^(-1/6)"
$ws.Range("E3:E66").Formula = "=B3^(-1/6)"
$ws.Range("E67:E130").Formula = "=B67^(-1/6)"
$ws.Range("E131:E194").Formula = "=B131^(-1/6)"
$ws.Range("E195:E207").Formula = "=B195^(-1/6)"

# Widen the new column and move the active selection.
$ws.Range("E:E").ColumnWidth = 17.2
$ws.Range("H30").Select() | Out-Null
